$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the sub-items of 7.a, 7.b, 7.c as passed ("x") -- i.e. "Added some more comments"
$passRows = 120, 121, 122, 123, 125, 126, 127, 128, 129, 131, 132, 133, 134, 135
foreach ($r in $passRows) {
    $ws.Range("C$r").Value = "x"
}

# Roll the pass mark up to the parent checklist rows now that all of their
# sub-items are complete.
$ws.Range("D119").Value = 2
$ws.Range("D124").Value = 2
$ws.Range("D130").Value = 2

# Update the view: zoom in further and move the selection/scroll position.
$excel.ActiveWindow.Zoom = 140
$ws.Range("B22").Select() | Out-Null
